$d = $word.ActiveDocument

# --- 1. Remove the double braces around the MMMM / YYYY placeholders in the
#        "date of order" line ("{{MMMM}}" -> "{MMMM}", "{{YYYY}}" -> "{YYYY}").
#        Each token is unique in the document, so a literal Find/Replace is safe.
$d.Content.Find.Execute("{{MMMM}}", $false, $false, $false, $false, $false, `
                         $true, 1, $false, "{MMMM}", 2) | Out-Null

$d.Content.Find.Execute("{{YYYY}}", $false, $false, $false, $false, $false, `
                         $true, 1, $false, "{YYYY}", 2) | Out-Null

# --- 2. Give the (single) table in the document explicit "no border" table
#        borders (top/left/bottom/right/insideH/insideV all set to "none").
$wdLineStyleNone    = 0
$wdColorAutomatic   = -16777216

$table = $d.Tables.Item(1)
$borders = $table.Borders

# Order matters for the underlying writer: set width/color first, then
# LineStyle last so the resulting border value is written as "none".
$borders.LineWidth = 0
$borders.Color = $wdColorAutomatic
$borders.LineStyle = $wdLineStyleNone
